# Listas sem duplicação de professores
# The teacher's weekly schedule listed a duplicated/encoded value
# "[-, 'MEC-1NA-Tec. Mat. Não Metal.', -, -]" in the E/F columns for the
# 20:00 and 20:50 rows (rows 18 and 19). Replace those values with a
# simple "-" to avoid showing the duplicated course/teacher list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "-"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "-"
